$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.784.86"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "2.774.58"
$ws.Range("E3").Value = "  +4.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.67"
$ws.Range("E5").Value = "  +5.09%  "
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  +4.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.577"
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.48"
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  +4.58%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.96"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.130"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.64"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "3.207.83"
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("D16").Value = "2.781.60"
$ws.Range("E16").Value = "  +4.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.880"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "51.720.23"
$ws.Range("E18").Value = "  +3.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.22"
$ws.Range("E19").Value = "  +9.36%  "
$ws.Range("E20").Value = "  +4.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.21"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "275.87"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.98"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.49"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.60"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.07"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.70"
$ws.Range("E33").Value = "  +3.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0818"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0382"
$ws.Range("E40").Value = "  +10.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("E41").Value = "  +24.36%  "
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.04"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.37"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "2.066.94"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.888"
$ws.Range("E50").Value = "  +12.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.87"
$ws.Range("E51").Value = "  -1.00%  "
